$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.507.75'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +2.59%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.513.58'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  +1.03%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  -0.03%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '590.86'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  +2.65%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '176.72'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +6.63%  '
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -0.02%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +2.30%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.511.53'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  +0.92%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.19%  '
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  +3.01%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '5.15'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +1.50%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.344'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.78%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '26.76'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +2.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.976.66'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.63%  '
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  +1.97%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.333.37'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +2.48%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.527.19'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +1.16%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '7.94'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +4.76%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.39'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +1.78%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '362.68'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +5.95%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.18'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +0.24%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.63'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +2.39%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '1.96'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +2.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '70.91'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +3.55%  '
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +3.26%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.00'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.40%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.644.62'
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '0.0₃0986'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.75%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '545.13'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +5.23%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.25'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.02%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  +3.15%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  +3.47%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.130'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.47%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.999'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -0.07%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.46'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.64%  '
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -1.16%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '18.75'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  +2.06%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '18.61'
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +1.82%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.356'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  +1.21%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '1.80'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '5.17'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +3.15%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.52'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +3.99%  '
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +0.01%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.560'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +2.01%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '146.49'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +0.36%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0₆0277'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +1.28%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '3.71'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  +1.09%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.68'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  +0.03%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.0754'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +0.59%  '
